$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.263.26'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '1.919.11'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8144'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.05'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3255'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.96'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07195'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7905'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08103'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').Value = '1.951.19'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.418'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.92'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '30.281.98'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.21'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.068'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '249.57'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007834'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.178.93'
$ws.Range('E21').Value = '  +0.79%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.162'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +19.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1676'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +18.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.489'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.91'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.02'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.178'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.383'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.61%  '
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.332'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05844'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.147'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.295'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7437'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.11%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.731'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9941'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01964'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.820'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4533'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.97'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.971'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8557'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.924'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.20'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.67%  '
$ws.Range('D48').Value = '1.022.68'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.983'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.95%  '
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.106'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +11.20%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.617'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.40%  '
